$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test")

# Add a new cell A7 storing the text value "234", using the same
# text number format as A5 (numFmtId 49 / "@") so it round-trips as
# a shared string rather than a number.
$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Value = "234"

# Update selection to mirror what Excel would leave selected after
# entering data in A7 (moves on to A8).
$ws.Range("A8").Select()
